$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New test-case rows (30-34): "interactive scenario with java client" ---
# Category column (H) filled first across the new rows - all share "CLDRJar"
$ws.Cells.Item(30,8).Value = "CLDRJar"
$ws.Cells.Item(31,8).Value = "CLDRJar"
$ws.Cells.Item(32,8).Value = "CLDRJar"
$ws.Cells.Item(33,8).Value = "CLDRJar"
$ws.Cells.Item(34,8).Value = "CLDRJar"

# Row 30 - number
$ws.Cells.Item(30,1).Value = 29
$ws.Cells.Item(30,2).Value = "Check the number data can be fetched"
$ws.Cells.Item(30,3).Value = "number"
$ws.Cells.Item(30,7).Value = 201703.54199999999
$ws.Cells.Item(30,7).NumberFormat = "#,##0.00"

# Row 32's expected value (plural) entered next
$ws.Cells.Item(32,7).Value = 'Il y a 100 000 fichiers sur "MyDisk".'

# Row 31 - percent
$ws.Cells.Item(31,1).Value = 30
$ws.Cells.Item(31,2).Value = "Check the percent data can be fetched"

# Row 32 - plural
$ws.Cells.Item(32,1).Value = 31
$ws.Cells.Item(32,2).Value = "Check the plurals data can be fetched"

$ws.Cells.Item(31,3).Value = "percent"
$ws.Cells.Item(31,7).Value = 0.23
$ws.Cells.Item(31,7).NumberFormat = "0%"

$ws.Cells.Item(32,3).Value = "plural"

# Row 33 - datetime
$ws.Cells.Item(33,1).Value = 32
$ws.Cells.Item(33,2).Value = "Check the datetime data can be fetched"
$ws.Cells.Item(33,3).Value = "datetime"

# Row 34 - currency
$ws.Cells.Item(34,1).Value = 33
$ws.Cells.Item(34,2).Value = "Check the currency data can be fetched"
$ws.Cells.Item(34,3).Value = "currency"

# Expected values for row 34 then row 33
$ws.Cells.Item(34,7).Value = 'US$201,703.54'
$ws.Cells.Item(33,7).Value = '2017年11月20日 GMT+8 下午1:39:24'

# --- AutoFilter over the original data range, with hidden filter-database name ---
$ws.Range("A1:H29").AutoFilter() | Out-Null
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$H`$29")
$fdb.Visible = $false

# --- View state: scrolled down, new cell selected ---
$ws.Range("A35").Select() | Out-Null
